# Removing totals in tags of measures in example files.
#
# For the K10+, K5 and SDQ sheets, the "tags" column (holding values like
# "total=25") is cleared out for every data row, leaving the column header
# itself ("k10p_tags" / "k5_tags" / "sdq_tags") intact. Clearing the cell
# contents (rather than just blanking the value) removes the now-unused
# <c> elements, and Excel automatically compacts the shared string table
# on save, dropping the "total=NN" strings that are no longer referenced
# anywhere in the workbook.

$wb = $excel.ActiveWorkbook

# K10+ sheet: column S (k10p_tags), data rows 2-6
$wsK10 = $wb.Worksheets.Item("K10+")
$wsK10.Range("S2:S6").ClearContents()

# K5 sheet: column J (k5_tags), data rows 2-11
$wsK5 = $wb.Worksheets.Item("K5")
$wsK5.Range("J2:J11").ClearContents()

# SDQ sheet: column BB (sdq_tags), data rows 2-22
$wsSDQ = $wb.Worksheets.Item("SDQ")
$wsSDQ.Range("BB2:BB22").ClearContents()
